# Issue#88: rename the "wt" and "dcin5" worksheets to clarify that they
# hold log2 expression data.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("wt").Name = "wt_log2_expression"
$wb.Worksheets.Item("dcin5").Name = "dcin5_log2_expression"
